# Work from 7/11-7/13 [1]
#
# The shared-strings header "Quote" is renamed to "Tool" (cell A2 of
# Sheet1), and the active window's view state is updated: the selection
# moves from C3 to A3 and the viewport scrolls so column G is the
# left-most visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "Quote" -> "Tool"
$ws.Range("A2").Value = "Tool"

# Update the view state: scroll so column G is left-most, then move the
# active selection to A3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
[void]$ws.Range("A3").Select()
